$d = $word.ActiveDocument

$replacements = @(
    @{old = "2024-02-11 Sunday"; new = "2024-02-12 Monday"},
    @{old = "471×9="; new = "703×3="},
    @{old = "182×4="; new = "222×3="},
    @{old = "676×8="; new = "800×2="},
    @{old = "396×5="; new = "144×7="},
    @{old = "401×5="; new = "152×9="},
    @{old = "169×2="; new = "840×4="},
    @{old = "303×2="; new = "430×2="},
    @{old = "971×4="; new = "644×4="},
    @{old = "113×4="; new = "183×3="},
    @{old = "665×8="; new = "636×4="},
    @{old = "885×8="; new = "239×5="},
    @{old = "600×2="; new = "295×5="},
    @{old = "441×3="; new = "677×9="},
    @{old = "877×7="; new = "821×6="},
    @{old = "179×4="; new = "759×9="},
    @{old = "145×4="; new = "365×8="},
    @{old = "845×8="; new = "772×8="},
    @{old = "944×5="; new = "602×6="},
    @{old = "743×8="; new = "742×2="},
    @{old = "180×7="; new = "270×3="},
    @{old = "565×5="; new = "891×8="},
    @{old = "215×9="; new = "886×5="},
    @{old = "461×6="; new = "156×2="},
    @{old = "558×6="; new = "522×9="},
    @{old = "343×6="; new = "149×4="}
)

foreach ($r in $replacements) {
    $d.Content.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
